$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Charges appartement: 2T 2017 (column C = "2T", row 6 = year 2017)
$ws.Range("C6").Value = 846.29

# Update the active cell selection to match the author's last cursor position
$ws.Range("C14").Select()
